$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing three data rows (5:7) - the table now only needs 3 data rows (2:4)
$ws.Rows("5:7").Delete()

# Row 2: FAPs / Hc / C5ar2 -> ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Hc"
$ws.Range("C2").Value = "C5ar2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.12672
$ws.Range("H2").Value = 0.38016
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.002689333333333334
$ws.Range("N2").Value = 0.008068000000000001
$ws.Range("O2").Value = 0.03638167388167388
$ws.Range("P2").Value = 0.03638167388167389
$ws.Range("Q2").Value = 0.00034079232
$ws.Range("R2").Value = 0.00306713088
$ws.Range("S2").Value = 0.03638167388167388
$ws.Range("T2").Value = 0.03638167388167389

# Row 3: FAPs / Hc / C5ar2 -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Hc"
$ws.Range("C3").Value = "C5ar2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.12672
$ws.Range("H3").Value = 0.38016
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.06844499999999999
$ws.Range("N3").Value = 0.205335
$ws.Range("O3").Value = 0.9259334415584415
$ws.Range("P3").Value = 0.9259334415584416
$ws.Range("Q3").Value = 0.008673350399999998
$ws.Range("R3").Value = 0.0780601536
$ws.Range("S3").Value = 0.9259334415584415
$ws.Range("T3").Value = 0.9259334415584416

# Row 4: FAPs / Hc / C5ar2 -> MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Hc"
$ws.Range("C4").Value = "C5ar2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.12672
$ws.Range("H4").Value = 0.38016
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.002785666666666667
$ws.Range("N4").Value = 0.008357
$ws.Range("O4").Value = 0.03768488455988456
$ws.Range("P4").Value = 0.03768488455988456
$ws.Range("Q4").Value = 0.00035299968
$ws.Range("R4").Value = 0.00317699712
$ws.Range("S4").Value = 0.03768488455988456
$ws.Range("T4").Value = 0.03768488455988456
